$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a running time-series log where a header block is
# repeated roughly every 36 rows, followed by one data row per State/UT.
# Append the next block: a header row (copied from the most recent header
# block at row 937 to keep the same bold/bordered formatting) followed by
# 35 state/UT rows for 14-10-2020.

$headerRow = 973
$ws.Range("A937:H937").Copy($ws.Range("A" + $headerRow + ":H" + $headerRow))
$ws.Cells.Item($headerRow, 1).Value = "States/UT"
$ws.Cells.Item($headerRow, 2).Value = "Active Cases"
$ws.Cells.Item($headerRow, 3).Value = "Active Cases Since Yesterday"
$ws.Cells.Item($headerRow, 4).Value = "Recovered Cases"
$ws.Cells.Item($headerRow, 5).Value = "Recovered Cases Since Yesterday"
$ws.Cells.Item($headerRow, 6).Value = "Deceased Cases"
$ws.Cells.Item($headerRow, 7).Value = "Deceased Cases Since Yesterday"
$ws.Cells.Item($headerRow, 8).Value = "Date"

$data = @(
    @("Andaman and Nicobar Islands", 199, 1, 3782, 12, 55, 0, "14-10-2020"),
    @("Andhra Pradesh", 42855, -1128, 714427, 5715, 6291, 35, "14-10-2020"),
    @("Arunachal Pradesh", 2960, 20, 9573, 170, 28, 4, "14-10-2020"),
    @("Assam", 28897, 458, 167059, 1020, 830, 4, "14-10-2020"),
    @("Bihar", 10835, 166, 187059, 1148, 961, 6, "14-10-2020"),
    @("Chandigarh", 1127, -43, 12007, 109, 197, 5, "14-10-2020"),
    @("Chhattisgarh", 27208, -213, 119352, 2812, 1306, 20, "14-10-2020"),
    @("Dadra and Nagar Haveli and Daman and Diu", 85, -17, 3081, 17, 2, 0, "14-10-2020"),
    @("Delhi", 21490, 955, 286880, 2036, 5854, 45, "14-10-2020"),
    @("Goa", 4316, -149, 34252, 554, 514, 3, "14-10-2020"),
    @("Gujarat", 15187, -227, 134990, 1375, 3584, 10, "14-10-2020"),
    @("Haryana", 10319, -82, 132382, 1154, 1601, 9, "14-10-2020"),
    @("Himachal Pradesh", 2507, -130, 15001, 311, 254, 3, "14-10-2020"),
    @("Jammu and Kashmir", 9866, -126, 73502, 796, 1340, 7, "14-10-2020"),
    @("Jharkhand", 7617, -159, 85314, 853, 805, 7, "14-10-2020"),
    @("Karnataka", 113478, -2317, 602505, 10421, 10123, 87, "14-10-2020"),
    @("Kerala", 95493, 1020, 207357, 7723, 1046, 21, "14-10-2020"),
    @("Ladakh", 969, 8, 4205, 79, 64, 0, "14-10-2020"),
    @("Madhya Pradesh", 14661, -271, 132429, 1708, 2671, 26, "14-10-2020"),
    @("Maharashtra", 205884, -7021, 1297252, 15356, 40701, 187, "14-10-2020"),
    @("Manipur", 2867, 111, 10829, 122, 97, 4, "14-10-2020"),
    @("Meghalaya", 2367, -67, 5406, 133, 65, 1, "14-10-2020"),
    @("Mizoram", 119, -37, 2093, 47, 0, 0, "14-10-2020"),
    @("Nagaland", 1513, 104, 5831, 18, 19, 1, "14-10-2020"),
    @("Odisha", 22892, -538, 232988, 2796, 1057, 17, "14-10-2020"),
    @("Puducherry", 4572, -45, 26865, 310, 567, 2, "14-10-2020"),
    @("Punjab", 8212, -364, 113105, 1006, 3894, 34, "14-10-2020"),
    @("Rajasthan", 21924, 253, 139616, 1768, 1679, 14, "14-10-2020"),
    @("Sikkim", 344, -40, 3014, 89, 59, 2, "14-10-2020"),
    @("Tamil Nadu", 43239, -508, 612320, 5117, 10371, 57, "14-10-2020"),
    @("Telengana", 23728, -480, 191269, 1918, 1241, 8, "14-10-2020"),
    @("Tripura", 3500, -238, 25041, 418, 318, 1, "14-10-2020"),
    @("Uttarakhand", 6576, -400, 48283, 674, 782, 20, "14-10-2020"),
    @("Uttar Pradesh", 38082, -733, 397570, 3662, 6466, 28, "14-10-2020"),
    @("West Bengal", 30988, 384, 265288, 3185, 5744, 62, "14-10-2020")
)

$row = 974
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
    $ws.Cells.Item($row, 8).Value = $item[7]
    $row = $row + 1
}
